# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values on the zh-cn and de-de report sheets (handback status report refresh).

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-18 05:55:56"
$wsZhCn.Range("G2").Value = "2016-01-18 05:56:49"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-18 05:56:08"
$wsDeDe.Range("G2").Value = "2016-01-18 05:57:11"
